$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Food")

# Add the new food items in the order they were first typed (this determines
# the order new shared-string entries get created in the workbook).
$newItems = @(
    @("Pork", "Meat"),
    @("Chicken", "Meat"),
    @("Tuna", "Fish"),
    @("Lobster", "Fish"),
    @("Pasta", "Grain"),
    @("Rice", "Grain"),
    @("Onion", "Veggie"),
    @("Lettuce", "Veggie")
)

foreach ($item in $newItems) {
    $row = $lo.ListRows.Add()
    $row.Range.Item(1, 1).Value = $item[0]
    $row.Range.Item(1, 2).Value = $item[1]
    $row.Range.Item(1, 3).Value = 10
    $row.Range.Item(1, 4).Value = 100
    $row.Range.Item(1, 5).Value = 1
}

# Sort the whole table by the Type column (ascending), like the author did.
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($lo.ListColumns.Item("Type").Range)
$lo.Sort.Header = 1
$lo.Sort.Apply()

# The validation ranges need to grow from row 7 to row 15 along with the table.
$ws.Range("C4:C7").Validation.Delete()
$ws.Range("F4:F7").Validation.Delete()
$ws.Range("E4:E7").Validation.Delete()

$ws.Range("C4:C15").Validation.Add(3, 1, 1, 'INDIRECT("ItemTypes[Type]")')
$ws.Range("F4:F15").Validation.Add(1, 1, 1, "1", "2")
$ws.Range("E4:E15").Validation.Add(1, 1, 1, "1", "100")

# Mirror the author's final selection.
$ws.Range("C3").Select()
